$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new column M, inheriting column L's formatting (row styles, borders,
# number formats, etc.) so every pre-existing row keeps its correct style index.
$ws.Columns("L:L").Copy() | Out-Null
$ws.Columns("M:M").Insert(-4161, 0)   # xlShiftToRight, xlFormatFromLeftOrAbove
$excel.CutCopyMode = $false

# Populate the new 2021 column with its own values (overwriting the
# duplicated-from-L placeholder values the Insert left behind).
$ws.Range("M4").Value = 2021
$ws.Range("M5").Value = 98
$ws.Range("M6").Value = 97
$ws.Range("M7").Value = 96

# Reset the view back to the top-left corner / default selection.
$ws.Range("A1").Select()
